$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed some bugs in CalcAdjacentPay: rows were re-sorted/recalculated,
# write the corrected (re-ordered) data rows back to the sheet.

$data = @(
    @(1001, 18, 30, 75, 60, 72),
    @(401,  9,  48, 67, 75, 45),
    @(101,  9,  30, 15, 60, 15),
    @(701,  3,  90, 45, 97, 15),
    @(1201, 2,  10, 10, 10, 10),
    @(1202, 2,  10, 10, 10, 10),
    @(901,  16, 15, 45, 60, 60),
    @(601,  9,  60, 67, 60, 42),
    @(1203, 3,  15, 15, 15, 15),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6,  45, 30, 60, 45),
    @(501,  9,  52, 30, 75, 45),
    @(201,  9,  30, 15, 45, 30),
    @(801,  3,  67, 65, 52, 45),
    @(502,  0,  4,  0,  0,  0),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0,  15, 30, 30, 0),
    @(1,    0,  2,  2,  2,  2),
    @(2,    0,  2,  2,  2,  2),
    @(802,  0,  4,  5,  4,  0),
    @(402,  0,  0,  4,  0,  0),
    @(602,  0,  0,  4,  0,  9),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
